# Auto-generated edit script applying the crypto price/volume refresh described in the commit diff.
# All target cells are plain text (inline strings) in the source sheet, so force a Text
# number format before assigning the value -- this stops Excel from re-interpreting
# numeric-looking strings (e.g. "1.00", "0.320", "5.80") as numbers and dropping the
# formatted trailing zeros / re-formatting punctuation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.706.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.515.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.53"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.514.23"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.360"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.975.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.426.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.91"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.524.34"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.70"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.97"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.643.04"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0897"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "464.36"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.86"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.56"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.320"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.72"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.29"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.12"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.22"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.99"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0738"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.80"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.34%  "
